# Add new tool "iProver v3.1 (CASC-27)" / run set "qbf-mode.SAT" results
# to the "runs" sheet (columns N:Q) and the "realtime" sheet (columns F:G).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "runs": new columns N (status), O (cputime), P (walltime), Q (memory)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("runs")

$sheet1Data = @(
  @(4, "sat", 38.64773622, 38.6714157126843, 120.844287999999),
  @(5, "unsat", 564.518411558, 564.137317843735, 66.3347199999999),
  @(6, "sat", 199.269575388, 199.181210082024, 112.336896),
  @(7, "TIMEOUT", 901.689430405, 901.064856108278, 100.306944),
  @(8, "TIMEOUT", 901.684762761, 901.065355155617, 46.534656),
  @(9, "TIMEOUT", 901.689481875, 901.033633742481, 29.925376),
  @(10, "TIMEOUT", 901.692693552, 901.04896248877, 32.018432),
  @(11, "TIMEOUT", 901.597290992, 901.041262026876, 63.893504),
  @(12, "TIMEOUT", 901.694522245, 901.05768045783, 42.356736),
  @(13, "TIMEOUT", 901.688784754, 901.049368847161, 153.755647999999),
  @(14, "TIMEOUT", 901.687665152, 901.049442827701, 48.3368959999999),
  @(15, "sat", 571.718645745, 571.321720331907, 139.038719999999),
  @(16, "TIMEOUT", 901.697626817, 901.061580371111, 87.4782719999999),
  @(17, "TIMEOUT", 901.595525803, 901.061649922281, 93.892608),
  @(18, "TIMEOUT", 901.691268091, 901.050942607224, 30.4455679999999),
  @(19, "TIMEOUT", 901.690128723, 901.069390814751, 32.526336),
  @(20, "TIMEOUT", 901.693456652, 901.049418501555, 35.0085119999999),
  @(21, "TIMEOUT", 901.693982873, 901.045751396566, 37.490688),
  @(22, "TIMEOUT", 901.694260824, 901.045493014156, 39.698432),
  @(23, "TIMEOUT", 901.587604041, 901.063851829618, 43.0858239999999),
  @(24, "TIMEOUT", 901.687648013, 901.057800382375, 36.2373119999999),
  @(25, "TIMEOUT", 901.674304308, 901.081781286746, 59.908096)
)

$ws1.Range("N1:Q1").Value2 = "iProver v3.1 (CASC-27)"
$ws1.Range("N2:Q2").Value2 = "qbf-mode.SAT"
$ws1.Range("N3").Value2 = "status"
$ws1.Range("O3").Value2 = "cputime (s)"
$ws1.Range("P3").Value2 = "walltime (s)"
$ws1.Range("Q3").Value2 = "memory (MB)"

foreach ($row in $sheet1Data) {
    $r = $row[0]
    $ws1.Cells.Item($r, 14).Value2 = $row[1]
    $ws1.Cells.Item($r, 15).Value2 = $row[2]
    $ws1.Cells.Item($r, 16).Value2 = $row[3]
    $ws1.Cells.Item($r, 17).Value2 = $row[4]
}

# Restore the selection shown in the sheet after the edit
[void]$ws1.Range("O1:O3").Select()

# ---------------------------------------------------------------
# Sheet "realtime": new columns F (status), G (cputime = runs!C + runs!O)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("realtime")

$sheet2Data = @(
  @(4, "sat"),
  @(5, "unsat"),
  @(6, "sat"),
  @(7, "TIMEOUT"),
  @(8, "TIMEOUT"),
  @(9, "TIMEOUT"),
  @(10, "TIMEOUT"),
  @(11, "TIMEOUT"),
  @(12, "TIMEOUT"),
  @(13, "TIMEOUT"),
  @(14, "TIMEOUT"),
  @(15, "sat"),
  @(16, "TIMEOUT"),
  @(17, "TIMEOUT"),
  @(18, "TIMEOUT"),
  @(19, "TIMEOUT"),
  @(20, "TIMEOUT"),
  @(21, "TIMEOUT"),
  @(22, "TIMEOUT"),
  @(23, "TIMEOUT"),
  @(24, "TIMEOUT"),
  @(25, "TIMEOUT")
)

$ws2.Range("F1:G1").Value2 = "iProver v3.1 (CASC-27)"
$ws2.Range("F2:G2").Value2 = "qbf-mode.SAT"
$ws2.Range("F3").Value2 = "status"
$ws2.Range("G3").Value2 = "cputime (s)"

foreach ($row in $sheet2Data) {
    $r = $row[0]
    $ws2.Cells.Item($r, 6).Value2 = $row[1]
    $ws2.Cells.Item($r, 7).Formula = "=runs!`$C$r+runs!`$O$r"
}

# Totals rows (27 solved-count, 28 uniquely-solved placeholder)
$ws2.Range("D27").Formula = "=COUNTIF(D4:D25,`"sat`") + COUNTIF(D4:D25,`"unsat`")"
$ws2.Range("F27").Formula = "=COUNTIF(F4:F25,`"sat`") + COUNTIF(F4:F25,`"unsat`")"
$ws2.Range("D28").Value2 = 0
$ws2.Range("F28").Value2 = 0

# Restore the selection shown in the sheet after the edit
[void]$ws2.Range("D29").Select()
